$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("parameters")

# --- Insert a new "requirements" column (D) ---
$ws.Columns.Item(4).Insert()

# --- Header row (row 1) ---
$ws.Range("D1").Value = "requirements"

# --- Row 2: mode ---
$ws.Range("B2").Value = "single_run"
$ws.Range("C2").Value = "instance_generation, single_run, combination_run"

# --- Rows 3-8: required for single_run and combination_run modes ---
$ws.Range("C3").ClearContents()
$ws.Range("C4").ClearContents()
$ws.Range("C5").ClearContents()
$ws.Range("B6").Value = 4
$ws.Range("C6").ClearContents()
$ws.Range("B7").Value = 120
$ws.Range("C7").ClearContents()
$ws.Range("C8").ClearContents()
$ws.Range("D3:D8").Value = "required for single_run and combination_run modes"

# --- Rows 10-20: required for instance_generation mode ---
$ws.Range("C10").ClearContents()
$ws.Range("C11").Value = "defines region size (e.g 7 sets  7x7 region)"
$ws.Range("C12").Value = "?"
$ws.Range("C13").Value = "??? Should we ignore this?" + [char]10 + "Do we have control on the number of water resources?"
$ws.Range("C14").Value = "0: no fire proof nodes, 1: add fire proof nodes" + [char]10 + "do we have control on the number of blocks? "
$ws.Range("C15").Value = "in units/sq km, 1: 0; 2: <6; 3: 6-50; 4: 50-741; 5: > 741"
$ws.Range("C16").Value = "0 if <50% vegetated, 1 if >50% vegetated"
$ws.Range("D10:D20").Value = "required for instance_generation mode"

# --- Formatting ---
$ws.Range("A1:D1").HorizontalAlignment = -4131
$ws.Range("A2:D8").HorizontalAlignment = -4131
$ws.Range("A9:I9").HorizontalAlignment = -4131
$ws.Range("A10:D20").HorizontalAlignment = -4131
$ws.Range("C13:C14").WrapText = $true
$ws.Rows.Item(13).RowHeight = 28.8
$ws.Rows.Item(14).RowHeight = 28.8

$ws.Range("B3").Select()
